# TC03_Canine_Filter_SampleType-NormalTissue.xlsx
# Replace the old "dbExcel" Neo4j query (shared by CasesTab/SamplesTab/FilesTab
# rows in column C) with the new consolidated filter query that also returns
# aliquot counts. Column B (StatQuery) text is unchanged; it only shifts its
# shared-string index because the old query string is dropped from the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @"
 MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])
        AND (s.study_disposition = 'Unrestricted')
        AND (size([]) = 0 OR s.clinical_study_type IN [])
        AND (size([]) = 0 OR demo.breed IN [])
        AND (size([]) = 0 OR demo.sex IN [])
        AND (size([]) = 0 OR demo.neutered_indicator IN [])
        AND (size([]) = 0 OR diag.disease_term IN [])
        AND (size([]) = 0 OR diag.primary_disease_site IN [])
        AND (size([]) = 0 OR diag.stage_of_disease IN [])
        AND (size([]) = 0 OR diag.best_response IN [])
    OPTIONAL MATCH (c)-->(co:cohort)
    OPTIONAL MATCH (f:file)-[*]->(c)
    OPTIONAL MATCH (f)-->(parent)
    OPTIONAL MATCH (samp:sample)-->(c)
    OPTIONAL MATCH (samp)<--(al:aliquot)
    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al
      WHERE (size(['Normal Tissue']) = 0 OR samp.summarized_sample_type IN ['Normal Tissue'])
        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])
        AND (size([]) = 0 OR samp.sample_site IN [])
        AND (size([]) = 0 OR head(labels(parent)) IN [])
        AND (size([]) = 0 OR f.file_type IN [])
        AND (size([]) = 0 OR f.file_format IN [])
    WITH c.case_id AS case_id,
         s.clinical_study_designation AS study_code,
         s.clinical_study_type AS study_type,
         co.cohort_description AS cohort,
         demo.breed AS breed,
         diag.disease_term AS diagnosis,
         diag.stage_of_disease AS stage_of_disease,
         diag.primary_disease_site AS disease_site,
         demo.patient_age_at_enrollment AS age,
         demo.sex AS sex,
         demo.neutered_indicator AS neutered_status,
         demo.weight AS weight,
         diag.best_response AS response_to_treatment,
         samp.sample_id AS sample_id,
         f.uuid AS file_id,
         al
    RETURN
COUNT(DISTINCT file_id) as number_of_files,
COUNT(DISTINCT sample_id) as number_of_sample,
COUNT(DISTINCT case_id) as number_of_cases,
COUNT(DISTINCT study_code) as number_of_study,
COUNT(DISTINCT al) as number_of_aliquot
    
"@

# dbExcel column (C) for CasesTab/SamplesTab/FilesTab all switch to the same
# new query text (replacing the old shared query, which disappears entirely
# since nothing references it anymore).
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# The wrapped rows now hold much longer text; Excel clamps to its maximum
# row height (409.6pt) instead of the old 244.8pt autofit height.
$ws.Rows(2).RowHeight = 409.6
$ws.Rows(3).RowHeight = 409.6
$ws.Rows(4).RowHeight = 409.6

# Scroll / selection moved down one row (B2 -> B3) and the view's top-left
# cell moved to A3.
[void]$ws.Activate()
[void]$ws.Range("B3").Select()
try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

# Best-effort restore of the saved window geometry recorded in the workbook.
try {
    $excel.ActiveWindow.Left = -108
    $excel.ActiveWindow.Top = -108
    $excel.ActiveWindow.Width = 23256
    $excel.ActiveWindow.Height = 12576
} catch {
}
